$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H70").Value = 8293.375
$ws.Range("I70").Value = 6199.8
$ws.Range("K70").Value = 18599.4
$ws.Range("M70").Value = -18329.4
$ws.Range("H73").Value = 8293.375
$ws.Range("I73").Value = 6199.8
$ws.Range("K73").Value = 18599.4
$ws.Range("M73").Value = -17663.4
$ws.Range("H86").Value = 7494.5
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 5250.8
$ws.Range("I88").Value = 2675
$ws.Range("J88").Value = 6968
$ws.Range("K88").Value = 2675
$ws.Range("L88").Value = 6968
$ws.Range("M88").Value = -2269
$ws.Range("N88").Value = -7780
$ws.Range("H89").Value = 7494.5
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 5250.8
$ws.Range("I91").Value = 2675
$ws.Range("J91").Value = 6968
$ws.Range("K91").Value = 2675
$ws.Range("L91").Value = 6968
$ws.Range("M91").Value = -1271
$ws.Range("N91").Value = -9776

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 39800
$ws.Range("J44").Value = 39800
$ws.Range("L44").Value = 39800
$ws.Range("N44").Value = -40776
$ws.Range("H98").Value = 37319.25
$ws.Range("J98").Value = 37319.25
$ws.Range("L98").Value = 37319.25
$ws.Range("N98").Value = -43309.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 500.75
$ws.Range("I22").Value = 200.8
$ws.Range("K22").Value = 200.8
$ws.Range("M22").Value = 149.2
$ws.Range("H31").Value = 5546
$ws.Range("I31").Value = 4019.7856
$ws.Range("J31").Value = 7682.7
$ws.Range("K31").Value = 4019.7856
$ws.Range("L31").Value = 7682.7
$ws.Range("M31").Value = -3724.7856
$ws.Range("N31").Value = -8272.700000000001
$ws.Range("H34").Value = 5546
$ws.Range("I34").Value = 4019.7856
$ws.Range("J34").Value = 7682.7
$ws.Range("K34").Value = 4019.7856
$ws.Range("L34").Value = 7682.7
$ws.Range("M34").Value = -3817.7856
$ws.Range("N34").Value = -8086.7
$ws.Range("H59").Value = 13049.5
$ws.Range("I59").Value = 13049.5
$ws.Range("K59").Value = 13049.5
$ws.Range("M59").Value = -11904.5
$ws.Range("H62").Value = 5420
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 6025
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 6025
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -7273
$ws.Range("H65").Value = 5420
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 6025
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 30125
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -36365
$ws.Range("H88").Value = 16111.25
$ws.Range("J88").Value = 16111.25
$ws.Range("L88").Value = 16111.25
$ws.Range("N88").Value = -16923.25
$ws.Range("H91").Value = 16111.25
$ws.Range("J91").Value = 16111.25
$ws.Range("L91").Value = 16111.25
$ws.Range("N91").Value = -18919.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 27000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H80").Value = 2916.3333
$ws.Range("I80").Value = 2916.3333
$ws.Range("K80").Value = 2916.3333
$ws.Range("M80").Value = -1918.3333
$ws.Range("H83").Value = 2916.3333
$ws.Range("I83").Value = 2916.3333
$ws.Range("K83").Value = 14581.6665
$ws.Range("M83").Value = -9589.666499999999
$ws.Range("H102").Value = 2594
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6897.0586
$ws.Range("I22").Value = 7173.077
$ws.Range("K22").Value = 7173.077
$ws.Range("M22").Value = -6878.077
$ws.Range("H27").Value = 6897.0586
$ws.Range("I27").Value = 7173.077
$ws.Range("K27").Value = 7173.077
$ws.Range("M27").Value = -7066.077
$ws.Range("H46").Value = 5333.25
$ws.Range("J46").Value = 5221.6523
$ws.Range("L46").Value = 5221.6523
$ws.Range("N46").Value = -5597.6523
$ws.Range("H55").Value = 777.2727
$ws.Range("J55").Value = 316.66666
$ws.Range("L55").Value = 316.66666
$ws.Range("N55").Value = -662.66666
$ws.Range("H82").Value = 2092.1538
$ws.Range("I82").Value = 1599.625
$ws.Range("K82").Value = 1599.625
$ws.Range("M82").Value = -1238.625
$ws.Range("H85").Value = 2092.1538
$ws.Range("I85").Value = 1599.625
$ws.Range("K85").Value = 1599.625
$ws.Range("M85").Value = -351.625
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H104").Value = 29476
$ws.Range("J104").Value = 29476
$ws.Range("L104").Value = 29476
$ws.Range("N104").Value = -36464

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H81").Value = 6199.875
$ws.Range("I81").Value = 1999.6666
$ws.Range("J81").Value = 8720
$ws.Range("K81").Value = 3999.3332
$ws.Range("L81").Value = 17440
$ws.Range("M81").Value = -2938.3332
$ws.Range("N81").Value = -19562
$ws.Range("H84").Value = 6199.875
$ws.Range("I84").Value = 1999.6666
$ws.Range("J84").Value = 8720
$ws.Range("K84").Value = 19996.666
$ws.Range("L84").Value = 87200
$ws.Range("M84").Value = -14692.666
$ws.Range("N84").Value = -97808
$ws.Range("H97").Value = 50572
$ws.Range("J97").Value = 50572
$ws.Range("L97").Value = 50572
$ws.Range("N97").Value = -52554
$ws.Range("H98").Value = 28000
$ws.Range("J98").Value = 28000
$ws.Range("L98").Value = 28000
$ws.Range("N98").Value = -33990

Write-Host "Applied Cuchulainn_Profits scheduled-runner updates"
